$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column O: header + two data rows (student id field)
$ws.Range("O1").Value = "Mã số sinh viên"
$ws.Range("O2").Value = 20161234
$ws.Range("O3").Value = 20161235

# Match the active cell selection shown in the diff (O4 selected, no range)
$ws.Range("O4").Select()
